# Trade #52 closed at 2026-02-17 12:49:19 - unknown UNKNOWN +0.000%
#
# This script updates the "live_trading_results" workbook to record the
# newly-closed trade #52:
#   - Summary sheet: Total Trades and Win Rate % are refreshed.
#   - Strategy Status sheet: MarketMaking row's Trades and Win Rate % are refreshed.
#   - All Trades sheet: a new row (53) is appended describing trade #52.
#   - MarketMaking sheet: the same new row (53) is appended.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B6").Value = 52      # Total Trades
$summary.Range("B9").Value = 40.38   # Win Rate %

# ---------------------------------------------------------------------------
# 2) Strategy Status sheet (MarketMaking row, row 4)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D4").Value = 52       # Trades
$status.Range("G4").Value = 40.38    # Win Rate %

# ---------------------------------------------------------------------------
# 3) Helper: append the row-53 trade record (same data for both sheets that
#    keep a trade log: "All Trades" and "MarketMaking").
# ---------------------------------------------------------------------------
function Add-Trade52Row($ws) {
    $row = 53

    $ws.Range("A$row").Value = 52

    # Columns B and C hold plain text (a date-looking string and a
    # time-looking string). Force text format first so Excel does not
    # auto-convert them into date/time serial values.
    $ws.Range("B$row").NumberFormat = "@"
    $ws.Range("B$row").Value = "2026-02-17"
    $ws.Range("B$row").Style = "Normal"

    $ws.Range("C$row").NumberFormat = "@"
    $ws.Range("C$row").Value = "12:49:13"
    $ws.Range("C$row").Style = "Normal"

    $ws.Range("D$row").Value = "MarketMaking"
    $ws.Range("E$row").Value = "UP"
    $ws.Range("F$row").Value = 0.97
    $ws.Range("G$row").Value = 0.97
    $ws.Range("H$row").Value = "CLOSED"
    $ws.Range("I$row").Value = 0
    $ws.Range("J$row").Value = 0
    $ws.Range("K$row").Value = 100.16
    $ws.Range("L$row").Value = 0
    $ws.Range("M$row").Value = 0
    $ws.Range("N$row").Value = 0.6
    $ws.Range("O$row").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P$row").Value = "early_exit"
    $ws.Range("Q$row").Value = 0.13
}

# ---------------------------------------------------------------------------
# 4) All Trades sheet
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Add-Trade52Row $allTrades

# ---------------------------------------------------------------------------
# 5) MarketMaking sheet (mirrors All Trades)
# ---------------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-Trade52Row $marketMaking
